$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width ---
$ws.Columns.Item(2).ColumnWidth = 13.1640625

# --- Row 3 header row height ---
$ws.Rows.Item(3).RowHeight = 31

# --- Clear old "Til faktura" label cell (A12) ---
$ws.Range("A12").ClearContents()

# =========================================================================
# New pizza-order table, rows 11-22
# =========================================================================

# --- Row 11: header ---
$ws.Range("B11:C11").Merge()
$ws.Range("B11").Value2 = "Lillekat2022-02"
$ws.Range("D11").Value2 = "Pris"

# Header formatting: font + fill copied from an existing themed header cell
$hdrSrc = $ws.Range("A4")
$hdrSrc.Copy()
$ws.Range("B11:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B11:C11").HorizontalAlignment = -4108
$ws.Range("D11").HorizontalAlignment = -4108

# Header borders: full box around B11:D11 (and split at the B/C merge + C/D boundary)
$ws.Range("B11:D11").Borders.Item(7).LineStyle = 1
$ws.Range("B11:D11").Borders.Item(7).Weight = 2
$ws.Range("B11:D11").Borders.Item(8).LineStyle = 1
$ws.Range("B11:D11").Borders.Item(8).Weight = 2
$ws.Range("B11:D11").Borders.Item(9).LineStyle = 1
$ws.Range("B11:D11").Borders.Item(9).Weight = 2
$ws.Range("D11").Borders.Item(10).LineStyle = 1
$ws.Range("D11").Borders.Item(10).Weight = 2

# Re-center after paste (paste special may reset alignment)
$ws.Range("B11:C11").HorizontalAlignment = -4108
$ws.Range("D11").HorizontalAlignment = -14142

# --- Rows 12-21: data rows ---
$ws.Cells.Item(12, 2).Value2 = 1
$ws.Cells.Item(12, 3).Value2 = 1
$ws.Cells.Item(12, 4).Formula = "=95*B12"

$ws.Cells.Item(13, 2).Value2 = 1
$ws.Cells.Item(13, 3).Value2 = 2
$ws.Cells.Item(13, 4).Value2 = 140

$ws.Cells.Item(14, 2).Value2 = 1
$ws.Cells.Item(14, 3).Value2 = 4
$ws.Cells.Item(14, 4).Value2 = 130

$ws.Cells.Item(15, 2).Value2 = 2
$ws.Cells.Item(15, 3).Value2 = 5
$ws.Cells.Item(15, 4).Formula = "=130*B15"

$ws.Cells.Item(16, 2).Value2 = 1
$ws.Cells.Item(16, 3).Value2 = "6+CF"
$ws.Cells.Item(16, 4).Value2 = 130

$ws.Cells.Item(17, 2).Value2 = 2
$ws.Cells.Item(17, 3).Value2 = 10
$ws.Cells.Item(17, 4).Formula = "=105*B17"

$ws.Cells.Item(18, 2).Value2 = 1
$ws.Cells.Item(18, 3).Value2 = 12
$ws.Cells.Item(18, 4).Value2 = 105

$ws.Cells.Item(19, 2).Value2 = 2
$ws.Cells.Item(19, 3).Value2 = 13
$ws.Cells.Item(19, 4).Formula = "=110*B19"

$ws.Cells.Item(20, 2).Value2 = 2
$ws.Cells.Item(20, 3).Value2 = 24
$ws.Cells.Item(20, 4).Formula = "=130*B20"

$ws.Cells.Item(21, 2).Value2 = 1
$ws.Cells.Item(21, 3).Value2 = 29
$ws.Cells.Item(21, 4).Value2 = 130

# Column borders for the data block (rows 12-21): left edge on B, right edge on D
$ws.Range("B12:B21").Borders.Item(7).LineStyle = 1
$ws.Range("B12:B21").Borders.Item(7).Weight = 2
$ws.Range("D12:D21").Borders.Item(10).LineStyle = 1
$ws.Range("D12:D21").Borders.Item(10).Weight = 2

# Right-align the "6+CF" label in C16
$ws.Range("C16").HorizontalAlignment = -4152

# Bottom border closing the data block under row 21
$ws.Range("B21").Borders.Item(9).LineStyle = 1
$ws.Range("B21").Borders.Item(9).Weight = 2
$ws.Range("C21").Borders.Item(9).LineStyle = 1
$ws.Range("C21").Borders.Item(9).Weight = 2
$ws.Range("D21").Borders.Item(9).LineStyle = 1
$ws.Range("D21").Borders.Item(9).Weight = 2

# --- Row 22: totals ---
$ws.Range("A22").Value2 = "Sum"
$ws.Range("B22").Formula = "=SUM(B12:B21)"
$ws.Range("D22").Formula = "=SUM(D12:D21)"

$ws.Range("A22:D22").Borders.Item(8).LineStyle = 1
$ws.Range("A22:D22").Borders.Item(8).Weight = 2
$ws.Range("A22:D22").Borders.Item(9).LineStyle = -4119
$ws.Range("A22").Borders.Item(10).LineStyle = 1
$ws.Range("A22").Borders.Item(10).Weight = 2
$ws.Range("B22").Borders.Item(7).LineStyle = 1
$ws.Range("B22").Borders.Item(7).Weight = 2
$ws.Range("D22").Borders.Item(10).LineStyle = 1
$ws.Range("D22").Borders.Item(10).Weight = 2

# --- View state ---
$ws.Range("B25").Select()
